$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.8
$ws.Range("H3").Value = 0.6848739495798319
$ws.Range("I3").Value = 0.07178111587982833
$ws.Range("J3").Value = 0.7
$ws.Range("K3").Value = 77.40000000000001

$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 16
$ws.Range("S3").Value = 42
$ws.Range("T3").Value = 67
$ws.Range("U3").Value = 210
$ws.Range("V3").Value = 927
$ws.Range("W3").Value = 916
$ws.Range("X3").Value = 890
$ws.Range("Y3").Value = 865
$ws.Range("Z3").Value = 722

$ws.Range("AF3").Value = 0.994635
$ws.Range("AG3").Value = 0.982833
$ws.Range("AH3").Value = 0.954936
$ws.Range("AI3").Value = 0.928112
$ws.Range("AJ3").Value = 0.774678
